$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the appID string used in rows 2 and 4 ("passive.income.nadi.myfirstdrawermenuproject"
# -> "...project2"). Using Cells.Replace edits the shared string text in place so every cell
# referencing it (B2 and B4) picks up the new value together.
$ws.Cells.Replace("passive.income.nadi.myfirstdrawermenuproject", "passive.income.nadi.myfirstdrawermenuproject2") | Out-Null

# Row height tweaks
$ws.Rows.Item(2).RowHeight = 12.8
$ws.Rows.Item(4).RowHeight = 23.85

# New row 14: another "bitcoin" / "com.hamxa.shaynachim" entry, same formatting as row 13
$ws.Range("A13:B13").Copy($ws.Range("A14:B14")) | Out-Null
$ws.Range("A14").Value = "bitcoin"
$ws.Range("B14").Value = "com.hamxa.shaynachim"

# New row 15: "affiliate marketing" / "affiliate.marketing.guide", same formatting, taller row
$ws.Range("A13:B13").Copy($ws.Range("A15:B15")) | Out-Null
$ws.Range("A15").Value = "affiliate marketing"
$ws.Range("B15").Value = "affiliate.marketing.guide"
$ws.Rows.Item(15).RowHeight = 24

# Move the active selection down to B16, below the newly added data
$ws.Range("B16").Select() | Out-Null
